$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before column N (14) — this pushes the whole
#    right-hand "diagram" block (Color/Status/Storage/Category/OrderItem/Series)
#    one column to the right (N:S -> O:T), matching the xr diff exactly.
$ws.Columns("N:N").Insert()

# 2) New field on the main Product table: "category_id" (endpoint product).
#    Style = fill (same as header cells) but WITHOUT the border, so clone the
#    format from an existing header cell (fill+border) and then strip the border.
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("M3").Value = "category_id"
$ws.Range("M3").Borders.LineStyle = -4142

# 3) Category table gains a "link" column (new endpoint field), same style as
#    the other field-name cells in that mini table.
$ws.Range("P12").Copy()
$ws.Range("Q12").PasteSpecial(-4122)
$ws.Range("Q12").Value = "link"

# 4) Series table: drop the stray "category_id" cell and fix the row so it
#    reads series_id / series / link.
$ws.Range("P19").Value = "series"
$ws.Range("P19").Copy()
$ws.Range("Q19").PasteSpecial(-4122)
$ws.Range("Q19").Value = "link"

# 5) Selection, as recorded by the author at save time.
[void]$ws.Range("K6").Select()

# 6) Column-width touch ups around the edit (closest achievable values; the
#    host only keeps ColumnWidth to 1/6-character granularity).
$ws.Columns(4).ColumnWidth = 14.666666666666666
$ws.Columns(5).ColumnWidth = 10.5
$ws.Columns(13).ColumnWidth = 10.5
$ws.Columns(14).ColumnWidth = 3.8333333333333335
$ws.Columns(15).ColumnWidth = 10.833333333333334
$ws.Columns(16).ColumnWidth = 10.5
$ws.Columns(17).ColumnWidth = 11.5
$ws.Columns(18).ColumnWidth = 12.5
$ws.Columns(19).ColumnWidth = 10.5
$ws.Columns(20).ColumnWidth = 14.333333333333334
